# Update countries & provincias Spain
# Applies the refreshed COVID-19 snapshot values and reorders the two
# country pairs whose case counts crossed over (Costa Rica/Nepal and
# Aruba/Jordania), plus bumps the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Last updated" banner ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Septiembre de 2020 a las 22:41"

# --- Costa Rica overtakes Nepal: swap the country labels on rows 57/58 ---
$ws.Range("A57").Value = "Costa Rica"
$ws.Range("A58").Value = "Nepal"

# --- Aruba overtakes Jordania: swap the country labels on rows 139/140 ---
$ws.Range("A139").Value = "Aruba"
$ws.Range("A140").Value = "Jordania"

# --- Refreshed statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6479721
$ws.Range("C4").Value = 19471
$ws.Range("D4").Value = 3740250
$ws.Range("E4").Value = 2546034
$ws.Range("G4").Value = 189
$ws.Range("H4").Value = 193437

# Row 10 - Sudafrica
$ws.Range("B10").Value = 639362
$ws.Range("C10").Value = 845
$ws.Range("D10").Value = 566555
$ws.Range("E10").Value = 57803
$ws.Range("G10").Value = 115
$ws.Range("H10").Value = 15004

# Row 24 - Alemania
$ws.Range("B24").Value = 253625
$ws.Range("C24").Value = 1901
$ws.Range("E24").Value = 17220
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 9405

# Row 28 - Israel
$ws.Range("B28").Value = 133975
$ws.Range("C28").Value = 3331
$ws.Range("D28").Value = 105455
$ws.Range("E28").Value = 27494

# Row 29 - Canada
$ws.Range("B29").Value = 132142
$ws.Range("C29").Value = 247
$ws.Range("D29").Value = 116459
$ws.Range("E29").Value = 6537

# Row 57 - now Costa Rica
$ws.Range("B57").Value = 48780
$ws.Range("C57").Value = 833
$ws.Range("D57").Value = 19083
$ws.Range("E57").Value = 29187
$ws.Range("G57").Value = 32
$ws.Range("H57").Value = 510

# Row 58 - now Nepal
$ws.Range("B58").Value = 47236
$ws.Range("C58").Value = 979
$ws.Range("D58").Value = 30677
$ws.Range("E58").Value = 16259
$ws.Range("G58").Value = 11
$ws.Range("H58").Value = 300

# Row 62 - Suiza
$ws.Range("D62").Value = 37700
$ws.Range("E62").Value = 4878

# Row 119 - Mozambique
$ws.Range("B119").Value = 4557
$ws.Range("C119").Value = 113
$ws.Range("D119").Value = 2697
$ws.Range("E119").Value = 1833

# Row 139 - now Aruba
$ws.Range("B139").Value = 2482
$ws.Range("C139").Value = 33
$ws.Range("D139").Value = 1244
$ws.Range("E139").Value = 1223
$ws.Range("H139").Value = 15

# Row 140 - now Jordania
$ws.Range("B140").Value = 2478
$ws.Range("C140").Value = 67
$ws.Range("D140").Value = 1817
$ws.Range("E140").Value = 644
$ws.Range("G140").Value = 1
$ws.Range("H140").Value = 17

# Row 189 - Barbados
$ws.Range("B189").Value = 179
$ws.Range("C189").Value = 1
$ws.Range("E189").Value = 18
